$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1399
$ws.Range("I86").Value = 998.3333
$ws.Range("K86").Value = 998.3333
$ws.Range("M86").Value = 124.6667

$ws.Range("H89").Value = 1399
$ws.Range("I89").Value = 998.3333
$ws.Range("K89").Value = 4991.6665
$ws.Range("M89").Value = 624.3334999999997

$ws.Range("H112").Value = 7049
$ws.Range("J112").Value = 7049
$ws.Range("L112").Value = 21147
$ws.Range("N112").Value = -23363

$ws.Range("H137").Value = 3511.9375
$ws.Range("I137").Value = 1575.7778
$ws.Range("J137").Value = 6001.2856
$ws.Range("K137").Value = 4727.3334
$ws.Range("L137").Value = 18003.8568
$ws.Range("M137").Value = -2177.3334
$ws.Range("N137").Value = -23103.8568

$ws.Range("H138").Value = 2418.561
$ws.Range("J138").Value = 2564.4119
$ws.Range("L138").Value = 7693.2357
$ws.Range("N138").Value = -17973.2357

$ws.Range("H141").Value = 4797.9
$ws.Range("I141").Value = 3973.75
$ws.Range("J141").Value = 5347.3335
$ws.Range("K141").Value = 11921.25
$ws.Range("L141").Value = 16042.0005
$ws.Range("M141").Value = -6741.25
$ws.Range("N141").Value = -26402.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4116.25
$ws.Range("I32").Value = 2280.625
$ws.Range("J32").Value = 10235
$ws.Range("K32").Value = 2280.625
$ws.Range("L32").Value = 10235
$ws.Range("M32").Value = -1993.625
$ws.Range("N32").Value = -10809

$ws.Range("H45").Value = 1442.75
$ws.Range("I45").Value = 1070.2727
$ws.Range("K45").Value = 1070.2727
$ws.Range("M45").Value = -693.2727

$ws.Range("H61").Value = 5194.7896
$ws.Range("I61").Value = 2924.1538
$ws.Range("K61").Value = 2924.1538
$ws.Range("M61").Value = -2712.1538

$ws.Range("H63").Value = 7001.6665
$ws.Range("I63").Value = 7001.6665
$ws.Range("K63").Value = 7001.6665
$ws.Range("M63").Value = -6315.6665

$ws.Range("H66").Value = 7001.6665
$ws.Range("I66").Value = 7001.6665
$ws.Range("K66").Value = 35008.3325
$ws.Range("M66").Value = -31576.3325

$ws.Range("H74").Value = 861.7586
$ws.Range("I74").Value = 829.3333
$ws.Range("K74").Value = 829.3333
$ws.Range("M74").Value = 44.66669999999999

$ws.Range("H77").Value = 861.7586
$ws.Range("I77").Value = 829.3333
$ws.Range("K77").Value = 4146.6665
$ws.Range("M77").Value = 221.3334999999997

$ws.Range("H92").Value = 49883
$ws.Range("J92").Value = 49883
$ws.Range("L92").Value = 49883
$ws.Range("N92").Value = -54875

$ws.Range("H132").Value = 1356.7037
$ws.Range("I132").Value = 1081.3684
$ws.Range("K132").Value = 3244.1052
$ws.Range("M132").Value = -714.1052

$ws.Range("H136").Value = 5194.7896
$ws.Range("I136").Value = 2924.1538
$ws.Range("K136").Value = 8772.4614
$ws.Range("M136").Value = -6222.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25000
$ws.Range("I82").Value = 5000
$ws.Range("K82").Value = 5000
$ws.Range("M82").Value = -4617

$ws.Range("H85").Value = 25000
$ws.Range("I85").Value = 5000
$ws.Range("K85").Value = 5000
$ws.Range("M85").Value = -3674

$ws.Range("H134").Value = 12144.739
$ws.Range("I134").Value = 13095.588
$ws.Range("K134").Value = 39286.764
$ws.Range("M134").Value = -36751.764

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 726.1429000000001
$ws.Range("I16").Value = 653.36365
$ws.Range("J16").Value = 993
$ws.Range("K16").Value = 653.36365
$ws.Range("L16").Value = 993
$ws.Range("M16").Value = -366.36365
$ws.Range("N16").Value = -1567

$ws.Range("H31").Value = 2763.8125
$ws.Range("I31").Value = 1962.7916
$ws.Range("K31").Value = 1962.7916
$ws.Range("M31").Value = -1667.7916

$ws.Range("H34").Value = 2763.8125
$ws.Range("I34").Value = 1962.7916
$ws.Range("K34").Value = 1962.7916
$ws.Range("M34").Value = -1760.7916

$ws.Range("H58").Value = 1450670.8
$ws.Range("I58").Value = 2416509.5
$ws.Range("K58").Value = 2416509.5
$ws.Range("M58").Value = -2416306.5

$ws.Range("H107").Value = 451.15384
$ws.Range("I107").Value = 396.81818
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 396.81818
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1523.18182
$ws.Range("N107").Value = -4590

$ws.Range("H113").Value = 726.1429000000001
$ws.Range("I113").Value = 653.36365
$ws.Range("J113").Value = 993
$ws.Range("K113").Value = 653.36365
$ws.Range("L113").Value = 993
$ws.Range("M113").Value = 1516.63635
$ws.Range("N113").Value = -5333

$ws.Range("H122").Value = 1359.2258
$ws.Range("I122").Value = 1315.619
$ws.Range("K122").Value = 3946.857
$ws.Range("M122").Value = -1496.857

$ws.Range("H132").Value = 2136.4814
$ws.Range("I132").Value = 1020.26666
$ws.Range("K132").Value = 3060.79998
$ws.Range("M132").Value = -530.7999799999998

$ws.Range("H134").Value = 843.46155
$ws.Range("I134").Value = 836.64
$ws.Range("J134").Value = 1014
$ws.Range("K134").Value = 2509.92
$ws.Range("L134").Value = 3042
$ws.Range("M134").Value = 25.07999999999993
$ws.Range("N134").Value = -8112

$ws.Range("H136").Value = 1450670.8
$ws.Range("I136").Value = 2416509.5
$ws.Range("K136").Value = 7249528.5
$ws.Range("M136").Value = -7246978.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 100.818184
$ws.Range("I12").Value = 48
$ws.Range("J12").Value = 112.55556
$ws.Range("K12").Value = 144
$ws.Range("L12").Value = 337.66668
$ws.Range("M12").Value = 29
$ws.Range("N12").Value = -683.66668

$ws.Range("H50").Value = 90973780
$ws.Range("I50").Value = 175092
$ws.Range("J50").Value = 142858740
$ws.Range("K50").Value = 525276
$ws.Range("L50").Value = 428576220
$ws.Range("M50").Value = -524795
$ws.Range("N50").Value = -428577182

$ws.Range("H53").Value = 90973780
$ws.Range("I53").Value = 175092
$ws.Range("J53").Value = 142858740
$ws.Range("K53").Value = 525276
$ws.Range("L53").Value = 428576220
$ws.Range("M53").Value = -524795
$ws.Range("N53").Value = -428577182

$ws.Range("H131").Value = 792.0599999999999
$ws.Range("I131").Value = 558.4
$ws.Range("J131").Value = 804.3579
$ws.Range("K131").Value = 1675.2
$ws.Range("L131").Value = 2413.0737
$ws.Range("M131").Value = 3364.8
$ws.Range("N131").Value = -12493.0737

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2751.8235
$ws.Range("J102").Value = 2551.75
$ws.Range("L102").Value = 2551.75
$ws.Range("N102").Value = -5795.75

$ws.Range("H126").Value = 2830424
$ws.Range("I126").Value = 4632346.5
$ws.Range("J126").Value = 127540.5
$ws.Range("K126").Value = 13897039.5
$ws.Range("L126").Value = 382621.5
$ws.Range("M126").Value = -13894569.5
$ws.Range("N126").Value = -387561.5

$ws.Range("H132").Value = 2027918.4
$ws.Range("I132").Value = 3207008.8
$ws.Range("K132").Value = 9621026.399999999
$ws.Range("M132").Value = -9618496.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1634.325
$ws.Range("I132").Value = 1409.826
$ws.Range("J132").Value = 1938.0588
$ws.Range("K132").Value = 4229.478
$ws.Range("L132").Value = 5814.1764
$ws.Range("M132").Value = -1699.478
$ws.Range("N132").Value = -10874.1764

$ws.Range("H136").Value = 3749.5334
$ws.Range("I136").Value = 2055.375
$ws.Range("K136").Value = 6166.125
$ws.Range("M136").Value = -3616.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1295.3
$ws.Range("I113").Value = 1017.1667
$ws.Range("J113").Value = 1712.5
$ws.Range("K113").Value = 3051.5001
$ws.Range("L113").Value = 5137.5
$ws.Range("M113").Value = -881.5001000000002
$ws.Range("N113").Value = -9477.5

$ws.Range("H122").Value = 126508.53
$ws.Range("I122").Value = 135366.28
$ws.Range("K122").Value = 406098.84
$ws.Range("M122").Value = -403648.84

$ws.Range("H126").Value = 6975.4243
$ws.Range("I126").Value = 6427.1
$ws.Range("J126").Value = 7819
$ws.Range("K126").Value = 19281.3
$ws.Range("L126").Value = 23457
$ws.Range("M126").Value = -16811.3
$ws.Range("N126").Value = -28397

$ws.Range("H132").Value = 1841.4642
$ws.Range("I132").Value = 1507.409
$ws.Range("K132").Value = 4522.227000000001
$ws.Range("M132").Value = -1992.227000000001

$ws.Range("H136").Value = 16341274
$ws.Range("I136").Value = 23149086
$ws.Range("K136").Value = 69447258
$ws.Range("M136").Value = -69444708
